$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old = $ws.Range("F37:H44")
$cf = $old.FormatConditions.Item(1)
Write-Host ("cf formula=" + $cf.Formula1)

$new = $excel.Union($ws.Range("H3"), $ws.Range("H8:H9"))
$new = $excel.Union($new, $ws.Range("H11:H12"))
$new = $excel.Union($new, $ws.Range("H14:H15"))
$new = $excel.Union($new, $ws.Range("H19:H25"))
$new = $excel.Union($new, $ws.Range("H29"))
$new = $excel.Union($new, $ws.Range("H31"))
$new = $excel.Union($new, $ws.Range("H33"))
$new = $excel.Union($new, $ws.Range("F37:H45"))
$new = $excel.Union($new, $ws.Range("J41:J45"))
Write-Host ("new Address=" + $new.Address())

$cf.ModifyAppliesToRange($new)
Write-Host "modified"
